$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.222.52'
$ws.Range("E2").Value = '  +0.88%  '

$ws.Range("D3").Value = '1.690.54'
$ws.Range("E3").Value = '  +0.75%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '''216.42'
$ws.Range("E5").Value = '  +0.64%  '

$ws.Range("D6").Value = '''0.522'
$ws.Range("E6").Value = '  +0.73%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").Value = '''23.02'
$ws.Range("E8").Value = '  +13.24%  '

$ws.Range("E9").Value = '  +4.73%  '

$ws.Range("D10").Value = '''0.0628'
$ws.Range("E10").Value = '  +1.47%  '

$ws.Range("D11").Value = '''0.0891'
$ws.Range("E11").Value = '  +0.41%  '

$ws.Range("D12").Value = '1.930.77'
$ws.Range("E12").Value = '  +0.89%  '

$ws.Range("D13").Value = '1.695.32'
$ws.Range("E13").Value = '  +0.69%  '

$ws.Range("E14").Value = '  +2.65%  '

$ws.Range("E15").Value = '  +5.02%  '

$ws.Range("D16").Value = '''67.58'
$ws.Range("E16").Value = '  +2.83%  '

$ws.Range("D17").Value = '27.238.06'

$ws.Range("D18").Value = '''237.59'
$ws.Range("E18").Value = '  +0.89%  '

$ws.Range("E19").Value = '  -0.97%  '

$ws.Range("D20").Value = '0.0₃0746'
$ws.Range("E20").Value = '  +1.65%  '

$ws.Range("E21").Value = '  -0.12%  '

$ws.Range("D22").Value = '''4.58'
$ws.Range("E22").Value = '  +3.00%  '

$ws.Range("D23").Value = '''9.65'
$ws.Range("E23").Value = '  +5.00%  '

$ws.Range("E24").Value = '  -2.41%  '

$ws.Range("D25").Value = '''148.25'
$ws.Range("E25").Value = '  +1.28%  '

$ws.Range("D26").Value = '''7.32'
$ws.Range("E26").Value = '  +1.35%  '

$ws.Range("D27").Value = '''16.52'
$ws.Range("E27").Value = '  +2.75%  '

$ws.Range("E28").Value = '  +1.21%  '

$ws.Range("D29").Value = '''0.999'
$ws.Range("E29").Value = '  -0.26%  '

$ws.Range("E30").Value = '  +1.14%  '

$ws.Range("E31").Value = '  +1.13%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '''3.41'
$ws.Range("E32").Value = '  +2.72%  '

$ws.Range("B33").Value = 'Maker'
$ws.Range("C33").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D33").Value = '1.573.97'
$ws.Range("E33").Value = '  +6.49%  '

$ws.Range("E34").Value = '  +2.52%  '

$ws.Range("D35").Value = '''1.69'
$ws.Range("E35").Value = '  +0.46%  '

$ws.Range("D36").Value = '''0.955'
$ws.Range("E36").Value = '  +5.88%  '

$ws.Range("E37").Value = '  +3.91%  '

$ws.Range("E38").Value = '  -1.14%  '

$ws.Range("E39").Value = '  -0.14%  '

$ws.Range("E40").Value = '  +4.30%  '

$ws.Range("D41").Value = '''69.59'
$ws.Range("E41").Value = '  +3.11%  '

$ws.Range("E42").Value = '  -1.01%  '

$ws.Range("E43").Value = '  -0.02%  '

$ws.Range("E44").Value = '  -2.32%  '

$ws.Range("D45").Value = '1.839.52'
$ws.Range("E45").Value = '  +1.22%  '

$ws.Range("D46").Value = '''0.789'
$ws.Range("E46").Value = '  +0.88%  '

$ws.Range("D47").Value = '''91.11'
$ws.Range("E47").Value = '  +0.74%  '

$ws.Range("D48").Value = '''1.61'
$ws.Range("E48").Value = '  +6.03%  '

$ws.Range("D49").Value = '0.0₆0110'
$ws.Range("E49").Value = '  +3.27%  '

$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '''0.105'
$ws.Range("E50").Value = '  +3.10%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '''8.22'
$ws.Range("E51").Value = '  +6.86%  '
